# Daily attendance processing - 2025-10-29 04:28:32
# Reorders the "Recorded By" (column G) name lists for specific rows
# to match the updated attendance processing log order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Range('G2').Value = 'system, backup@backdoor.com, System'
$ws.Range('G3').Value = 'dnasr281@gmail.com, System'
$ws.Range('G4').Value = 'backup@backdoor.com, System'
$ws.Range('G6').Value = 'dnasr281@gmail.com, System'
$ws.Range('G10').Value = 'dnasr281@gmail.com, System'
$ws.Range('G11').Value = 'dnasr281@gmail.com, System'
$ws.Range('G12').Value = 'dnasr281@gmail.com, System'
$ws.Range('G13').Value = 'dnasr281@gmail.com, System'
$ws.Range('G14').Value = 'dnasr281@gmail.com, System'
$ws.Range('G15').Value = 'dnasr281@gmail.com, System'
$ws.Range('G17').Value = 'dnasr281@gmail.com, System'
$ws.Range('G18').Value = 'dnasr281@gmail.com, System'
$ws.Range('G19').Value = 'dnasr281@gmail.com, System'
$ws.Range('G20').Value = 'dnasr281@gmail.com, System'
$ws.Range('G21').Value = 'dnasr281@gmail.com, System'
$ws.Range('G22').Value = 'dnasr281@gmail.com, System'
$ws.Range('G24').Value = 'dnasr281@gmail.com, System'
$ws.Range('G29').Value = 'system, backup@backdoor.com, System'
$ws.Range('G30').Value = 'dnasr281@gmail.com, System'
$ws.Range('G31').Value = 'backup@backdoor.com, System'
$ws.Range('G33').Value = 'dnasr281@gmail.com, System'
$ws.Range('G37').Value = 'dnasr281@gmail.com, System'
$ws.Range('G38').Value = 'dnasr281@gmail.com, System'
$ws.Range('G39').Value = 'dnasr281@gmail.com, System'
$ws.Range('G40').Value = 'dnasr281@gmail.com, System'
$ws.Range('G41').Value = 'dnasr281@gmail.com, System'
$ws.Range('G42').Value = 'dnasr281@gmail.com, System'
$ws.Range('G44').Value = 'dnasr281@gmail.com, System'
$ws.Range('G45').Value = 'dnasr281@gmail.com, System'
$ws.Range('G46').Value = 'dnasr281@gmail.com, System'
$ws.Range('G47').Value = 'dnasr281@gmail.com, System'
$ws.Range('G48').Value = 'dnasr281@gmail.com, System'
$ws.Range('G49').Value = 'dnasr281@gmail.com, System'
$ws.Range('G51').Value = 'dnasr281@gmail.com, System'
$ws.Range('G56').Value = 'system, backup@backdoor.com, System'
$ws.Range('G57').Value = 'dnasr281@gmail.com, System'
$ws.Range('G58').Value = 'backup@backdoor.com, System'
$ws.Range('G60').Value = 'dnasr281@gmail.com, System'
$ws.Range('G64').Value = 'dnasr281@gmail.com, System'
$ws.Range('G65').Value = 'dnasr281@gmail.com, System'
$ws.Range('G66').Value = 'dnasr281@gmail.com, System'
$ws.Range('G67').Value = 'dnasr281@gmail.com, System'
$ws.Range('G68').Value = 'dnasr281@gmail.com, System'
$ws.Range('G69').Value = 'dnasr281@gmail.com, System'
$ws.Range('G71').Value = 'dnasr281@gmail.com, System'
$ws.Range('G72').Value = 'dnasr281@gmail.com, System'
$ws.Range('G73').Value = 'dnasr281@gmail.com, System'
$ws.Range('G74').Value = 'dnasr281@gmail.com, System'
$ws.Range('G75').Value = 'dnasr281@gmail.com, System'
$ws.Range('G76').Value = 'dnasr281@gmail.com, System'
$ws.Range('G78').Value = 'dnasr281@gmail.com, System'
$ws.Range('G86').Value = 'dnasr281@gmail.com, System'
$ws.Range('G87').Value = 'dnasr281@gmail.com, System'
$ws.Range('G88').Value = 'dnasr281@gmail.com, System'
$ws.Range('G89').Value = 'dnasr281@gmail.com, System'
$ws.Range('G90').Value = 'dnasr281@gmail.com, admin@admin.com'
$ws.Range('G93').Value = 'dnasr281@gmail.com, System'
$ws.Range('G95').Value = 'dnasr281@gmail.com, System'
$ws.Range('G96').Value = 'dnasr281@gmail.com, System'
$ws.Range('G97').Value = 'dnasr281@gmail.com, System'
$ws.Range('G99').Value = 'dnasr281@gmail.com, System'
$ws.Range('G102').Value = 'dnasr281@gmail.com, System'
$ws.Range('G112').Value = 'dnasr281@gmail.com, System'
$ws.Range('G113').Value = 'dnasr281@gmail.com, System'
$ws.Range('G114').Value = 'dnasr281@gmail.com, System'
$ws.Range('G115').Value = 'dnasr281@gmail.com, System'
$ws.Range('G116').Value = 'dnasr281@gmail.com, admin@admin.com'
$ws.Range('G119').Value = 'dnasr281@gmail.com, System'
$ws.Range('G121').Value = 'dnasr281@gmail.com, System'
$ws.Range('G122').Value = 'dnasr281@gmail.com, System'
$ws.Range('G123').Value = 'dnasr281@gmail.com, System'
$ws.Range('G125').Value = 'dnasr281@gmail.com, System'
$ws.Range('G128').Value = 'dnasr281@gmail.com, System'
$ws.Range('G138').Value = 'dnasr281@gmail.com, System'
$ws.Range('G139').Value = 'dnasr281@gmail.com, System'
$ws.Range('G140').Value = 'dnasr281@gmail.com, System'
$ws.Range('G141').Value = 'dnasr281@gmail.com, System'
$ws.Range('G142').Value = 'dnasr281@gmail.com, admin@admin.com'
$ws.Range('G145').Value = 'dnasr281@gmail.com, System'
$ws.Range('G147').Value = 'dnasr281@gmail.com, System'
$ws.Range('G148').Value = 'dnasr281@gmail.com, System'
$ws.Range('G149').Value = 'dnasr281@gmail.com, System'
$ws.Range('G151').Value = 'dnasr281@gmail.com, System'
$ws.Range('G154').Value = 'dnasr281@gmail.com, System'
